# Update "中国实际利用世界外商直接投资金额" sheet:
#   - drop the 2000-2009 rows
#   - keep 2010-2020 (shifted up to rows 2-12)
#   - append two new years, 2021 and 2022 (rows 13-14)
# End result: A1:B14, row 1 is the header, rows 2-14 are year/value pairs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$years  = @(2010, 2011, 2012, 2013, 2014, 2015, 2016, 2017, 2018, 2019, 2020, 2021, 2022)
$values = @(10573524, 11601100, 11171614, 11758260, 11956156, 12626660, 12600142, 13103513, 13496589, 13813462, 14436926, 17348331, 18913241)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = "$($years[$i])年"
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Remove the now-unused trailing rows (old 2013-2020 data that used to live
# at rows 15-22) so the sheet shrinks back down to A1:B14.
$ws.Rows("15:22").Delete()
